# "Inclusão dos links do github"
# Point the four "Grupo N" labels (2023.1 column) at their GitHub repos.
# Doing this through TextRange.ActionSettings(ppMouseClick).Hyperlink.Address
# makes PowerPoint allocate new hyperlink relationships (rId7..rId10) instead
# of overwriting the relationships (rId3..rId6) still used by the 2023.2
# column's identically-labelled textboxes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$links = @{
    "TextBox 30" = "https://github.com/ifpb-redes-2023-1/grupo-1"
    "TextBox 31" = "https://github.com/ifpb-redes-2023-1/grupo-2"
    "TextBox 32" = "https://github.com/ifpb-redes-2023-1/grupo-3"
    "TextBox 33" = "https://github.com/ifpb-redes-2023-1/grupo-4"
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($links.ContainsKey($shape.Name)) {
        $actionSetting = $shape.TextFrame.TextRange.ActionSettings.Item(1)
        $actionSetting.Hyperlink.Address = $links[$shape.Name]
    }
}
